$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 8 (the blank separator row),
# pushing the payment table, chart anchor, etc. down by two rows.
$ws.Rows("8:9").Insert()

# --- New row 8: "Pension Present Worth" ----------------------------------
# Copy the shaded label formatting from D4 ("(g)Growth Rate(%)") onto D8,
# then set its text.
$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "Pension Present Worth"

# Copy the shaded value formatting from E4 onto E8, then fix it up to be a
# left-aligned currency figure (matching the rest of the $ column) instead
# of a right-aligned percentage, and set its formula.
$ws.Range("E4").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").NumberFormat = $ws.Range("E7").NumberFormat
$ws.Range("E8").HorizontalAlignment = -4131
$ws.Range("E8").Formula = "=PV(E5,E6,-E3)"

# --- New row 9: "Extra Savings" -------------------------------------------
# Row-insert already copied the plain bold-label / currency formatting down
# from row 7 onto row 9, matching the "(P)Present Worth" row above it, so we
# only need to fill in the values.
$ws.Range("D9").Value = "Extra Savings"
$ws.Range("E9").Formula = "=E7-E8"

# --- Resize the "Table5" list object / named table to its new location ----
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A11:C37"))

# --- Re-point the chart's cached series ranges at the table's new rows ----
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,'Retirement Savings'!`$A`$12:`$A`$37,'Retirement Savings'!`$C`$12:`$C`$37,1)"

# --- Update the print area to cover the two extra rows ---------------------
$ws.PageSetup.PrintArea = "`$A`$1:`$E`$38"

# --- Selection cosmetic change seen in the saved file -----------------------
$ws.Range("G9").Select()
